$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ B=87343;  C=118090; D=147639; E=250292; F=432068; G=76071;  H=65617 }
    3  = @{ B=87724;  C=126958; D=149037; E=257182; F=433406; G=76270;  H=65875 }
    4  = @{ B=98681;  C=127218; D=167936; E=289474; F=487998; G=85711;  H=74305 }
    5  = @{ B=139291; C=199161; D=233279; E=390918; F=582079; G=121509; H=105689 }
    6  = @{ B=239505; C=319032; D=419898; E=585262; F=797553; G=209861; H=183156 }
    7  = @{ B=62635;  C=90275;  D=112092; E=211865; F=309103; G=59541;  H=53586 }
    8  = @{ B=63811;  C=91166;  D=112056; E=213610; F=314841; G=59256;  H=55306 }
    9  = @{ B=73237;  C=103198; D=129377; E=247841; F=363534; G=69278;  H=64893 }
    10 = @{ B=89257;  C=125473; D=154100; E=298385; F=425395; G=83228;  H=77960 }
    11 = @{ B=74241;  C=100375; D=125492; E=212748; F=367257; G=64661;  H=55774 }
}

foreach ($row in $values.Keys) {
    foreach ($col in $values[$row].Keys) {
        $ws.Range("$col$row").Value = $values[$row][$col]
    }
}

$ws.Range("D16").Select()
